$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.548.95"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.833.78"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("D4").Formula = "'0.9990"
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Formula = "'331.40"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Formula = "'0.9973"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").Formula = "'0.4479"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("D8").Formula = "'0.3799"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Formula = "'44.99"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Formula = "'0.07791"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Formula = "'22.37"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Formula = "'0.9966"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").Formula = "'6.377"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Formula = "'7.608"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "1.839.13"
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Formula = "'93.00"
$ws.Range("E17").Value = "  +14.50%  "
$ws.Range("D18").Formula = "'0.00001090"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Formula = "'0.06366"
$ws.Range("E19").Value = "  -5.64%  "
$ws.Range("D20").Formula = "'0.9972"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").Formula = "'17.71"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Formula = "'6.418"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").Formula = "'0.5415"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "28.597.65"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Formula = "'11.90"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").Formula = "'2.259"
$ws.Range("E26").Value = "  -7.69%  "
$ws.Range("D27").Formula = "'21.01"
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").Formula = "'154.52"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").Formula = "'2.385"
$ws.Range("E29").Value = "  +1.02%  "
$ws.Range("D30").Value = "2.047.02"
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("D31").Formula = "'129.89"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").Formula = "'1.218"
$ws.Range("E32").Value = "  -4.38%  "
$ws.Range("D33").Formula = "'5.892"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").Formula = "'0.09298"
$ws.Range("E34").Value = "  +0.94%  "
$ws.Range("D35").Formula = "'3.681"
$ws.Range("E35").Value = "  -7.28%  "
$ws.Range("D36").Formula = "'12.93"
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("D37").Formula = "'0.02374"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").Formula = "'0.2211"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Formula = "'0.6690"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Formula = "'5.233"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Formula = "'0.06295"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Formula = "'1.202"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Formula = "'8.187"
$ws.Range("E43").Value = "  +1.15%  "

# Row 44/45: swap Frax and WEMIXTOKEN rows with updated values
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Formula = "'1.408"
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Formula = "'0.9964"
$ws.Range("E45").Value = "  -0.34%  "

$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Formula = "'0.6163"
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("D48").Formula = "'3.789"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").Formula = "'128.30"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Formula = "'2.054"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Formula = "'79.76"
$ws.Range("E51").Value = "  +1.86%  "

Write-Host "Cryptos list updated"
